# Scénarii et Diagramme commit:
# "Les scénarii sont terminé" -> the "Scénarii - F3" test-scenario sheet is
# filled in: it used to contain four empty placeholder sub-sections, it now
# contains one finished sub-section ("Paramétrage des acquisitions" /
# "Modification") with its four completed steps, and the unused placeholder
# sub-sections are removed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Scénarii - F3": finish the scenario, drop the unused templates
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Scénarii - F3")
$ws3.Activate()

# Only the first sub-section (rows 4-7) is kept & completed; the three other
# empty placeholder sub-sections (old rows 8-19) are deleted altogether.
$ws3.Rows("9:19").Delete()

# The old second sub-section header (row 8, merged B8:E8) becomes the 4th
# step of the remaining sub-section instead: unmerge it and give it the same
# look as the other data rows (5-7).
$ws3.Range("B8:E8").UnMerge()
$ws3.Range("B5:E5").Copy()
$ws3.Range("B8").PasteSpecial(-4122)

# Fill in the scenario title, sub-section title and the 4 step descriptions.
$ws3.Range("B2").Value = "Paramétrage des acquisitions"
$ws3.Range("B4").Value = "Modification "
$ws3.Range("C5").Value = "Réception des choix par défaut du contrôleur "
$ws3.Range("C6").Value = "Permettre la modification de ceux-ci "
$ws3.Range("C7").Value = "Permettre la modification des paramètre de la carte "
$ws3.Range("B8").Value = 4
$ws3.Range("C8").Value = "Lecture des données "

$ws3.Range("C11").Select()

# ---------------------------------------------------------------------
# Sheet "Scénarii - F5": leave content untouched, just move the cursor
# (last sheet the author was working in when the file was saved).
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Scénarii - F5")
$ws5.Activate()
$ws5.Range("C37").Select()
